# Auto-generated edit script for EDV/TLT/ZROZ historical data update
$wb = $excel.ActiveWorkbook

# ---- Sheet: EDV ----
$ws = $wb.Worksheets.Item("EDV")

# Updated Adj Close / OHLCV values (refreshed from source)
$EDVUpdates = @{
    "F7" = 84.340622
    "F8" = 86.238838
    "F10" = 87.636009
    "F12" = 89.44750999999999
    "F13" = 88.77301
    "F15" = 86.238838
    "F17" = 88.22378500000001
    "F19" = 87.664917
    "F21" = 87.934715
    "F22" = 89.534233
    "F24" = 88.08889000000001
    "F28" = 85.602898
    "F34" = 83.155449
    "F36" = 82.19188699999999
    "F44" = 82.403862
    "F45" = 83.11689800000001
    "F47" = 83.290352
    "F50" = 85.082573
    "F51" = 87.154228
    "F53" = 87.99252300000001
    "F54" = 86.82662999999999
    "F60" = 85.082573
    "F62" = 85.381271
    "C202" = 65
    "D202" = 62.799999
    "E202" = 62.830002
    "F202" = 62.830002
    "G202" = 1440700
}
foreach ($key in $EDVUpdates.Keys) {
    $ws.Range($key).Value = $EDVUpdates[$key]
}

# New row 203
$ws.Range("A203").Value = "'2023-10-20"
$ws.Range("B203").Value = 62.950001
$ws.Range("C203").Value = 63.509998
$ws.Range("D203").Value = 62.66
$ws.Range("E203").Value = 63.16
$ws.Range("F203").Value = 63.16
$ws.Range("G203").Value = 786400

# New row 204
$ws.Range("A204").Value = "'2023-10-23"
$ws.Range("B204").Value = 62.950001
$ws.Range("C204").Value = 62.779999
$ws.Range("D204").Value = 62.5
$ws.Range("E204").Value = 62.5
$ws.Range("F204").Value = 62.5
$ws.Range("G204").Value = 109191

# ---- Sheet: TLT ----
$ws = $wb.Worksheets.Item("TLT")

# Updated Adj Close / OHLCV values (refreshed from source)
$TLTUpdates = @{
    "F2" = 99.01917299999999
    "F4" = 100.795403
    "F5" = 102.649696
    "F6" = 103.19622
    "F7" = 101.488319
    "F8" = 103.137657
    "F10" = 104.181923
    "F11" = 103.508522
    "F15" = 103.157188
    "F16" = 104.640617
    "F17" = 104.894371
    "F18" = 104.406395
    "F19" = 104.142883
    "F20" = 103.762283
    "F21" = 104.59182
    "F23" = 105.986115
    "F25" = 103.628044
    "F26" = 102.796364
    "F28" = 102.307144
    "F29" = 101.162346
    "F31" = 101.778763
    "F32" = 100.829681
    "F33" = 99.401123
    "F34" = 100.174103
    "F35" = 98.22699
    "F38" = 98.794487
    "F39" = 99.068443
    "F40" = 99.518547
    "F43" = 99.940315
    "F45" = 99.773567
    "F46" = 99.871651
    "F47" = 100.116867
    "F48" = 103.569511
    "F50" = 102.098213
    "F51" = 104.069763
    "F52" = 103.255638
    "F55" = 102.990791
    "F56" = 104.364014
    "F57" = 104.364014
    "F59" = 102.343422
    "F60" = 102.529793
    "F62" = 102.794617
    "F63" = 104.334595
    "F64" = 104.825287
    "F65" = 105.346466
    "F66" = 106.447815
    "F67" = 106.723145
    "F68" = 105.002281
    "F69" = 105.21862
    "F71" = 104.284447
    "F77" = 102.661903
    "F79" = 105.179276
    "F80" = 104.058266
    "F81" = 103.025749
    "F82" = 104.687607
    "F83" = 101.668991
    "F84" = 104.202576
    "F85" = 104.784233
    "F86" = 103.7491
    "F89" = 101.590134
    "F90" = 102.575958
    "F91" = 103.660378
    "F92" = 102.792839
    "F93" = 101.728142
    "F94" = 101.422531
    "F96" = 100.377556
    "F97" = 99.66773999999999
    "F98" = 99.312859
    "F99" = 99.598732
    "F100" = 99.10581999999999
    "F101" = 98.85936700000001
    "F104" = 101.530968
    "F108" = 101.217644
    "F109" = 99.71517900000001
    "F110" = 100.881554
    "F111" = 100.743172
    "F112" = 101.039719
    "F113" = 100.041374
    "F116" = 101.415321
    "F117" = 102.117119
    "F120" = 102.136894
    "F121" = 102.245636
    "F122" = 101.978745
    "F123" = 102.413673
    "F124" = 100.565254
    "F129" = 98.20117999999999
    "F130" = 98.330017
    "F132" = 99.935654
    "F138" = 100.797935
    "F139" = 100.827667
    "F141" = 100.272629
    "F145" = 99.162575
    "F151" = 96.096504
    "F152" = 96.59343699999999
    "F153" = 95.00324999999999
    "F160" = 91.95210299999999
    "F162" = 94.953568
    "F163" = 94.327438
    "F164" = 94.635536
    "F165" = 94.734909
    "G201" = 63724600
    "D202" = 82.739998
    "E202" = 82.769997
    "F202" = 82.769997
    "G202" = 87696900
}
foreach ($key in $TLTUpdates.Keys) {
    $ws.Range($key).Value = $TLTUpdates[$key]
}

# New row 203
$ws.Range("A203").Value = "'2023-10-20"
$ws.Range("B203").Value = 82.989998
$ws.Range("C203").Value = 83.540001
$ws.Range("D203").Value = 82.769997
$ws.Range("E203").Value = 83.239998
$ws.Range("F203").Value = 83.239998
$ws.Range("G203").Value = 52162600

# New row 204
$ws.Range("A204").Value = "'2023-10-23"
$ws.Range("B204").Value = 82.989998
$ws.Range("C204").Value = 82.894997
$ws.Range("D204").Value = 82.58000199999999
$ws.Range("E204").Value = 82.58429700000001
$ws.Range("F204").Value = 82.58429700000001
$ws.Range("G204").Value = 3530145

# ---- Sheet: ZROZ ----
$ws = $wb.Worksheets.Item("ZROZ")

# Updated Adj Close / OHLCV values (refreshed from source)
$ZROZUpdates = @{
    "G201" = 506700
    "C202" = 67.519997
    "D202" = 64.800003
    "E202" = 64.83000199999999
    "F202" = 64.83000199999999
    "G202" = 781500
}
foreach ($key in $ZROZUpdates.Keys) {
    $ws.Range($key).Value = $ZROZUpdates[$key]
}

# New row 203
$ws.Range("A203").Value = "'2023-10-20"
$ws.Range("B203").Value = 64.93000000000001
$ws.Range("C203").Value = 65.449997
$ws.Range("D203").Value = 64.540001
$ws.Range("E203").Value = 65.150002
$ws.Range("F203").Value = 65.150002
$ws.Range("G203").Value = 550100

# New row 204
$ws.Range("A204").Value = "'2023-10-23"
$ws.Range("B204").Value = 64.93000000000001
$ws.Range("C204").Value = 64.5
$ws.Range("D204").Value = 64.26840199999999
$ws.Range("E204").Value = 64.277496
$ws.Range("F204").Value = 64.277496
$ws.Range("G204").Value = 66718
